$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.507.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.434.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.12'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.512'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.504'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.14'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.42%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0800'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.70'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.94'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.816.93'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.444.38'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.839'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.484.69'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.42'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.79'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.75'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.14%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.23'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +4.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.15'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.38'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +10.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.120'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +14.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0763'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.88%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.90'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.60%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.51'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '128.66'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +17.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.90'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.75'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -9.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0291'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.957.17'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.18'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.90'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.49'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.67'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +9.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.30'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.67'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.62'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.74%  '
